$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns store plain text values (e.g. "26.431.87",
# "1.000", "0.000007159") that Excel would otherwise auto-convert/normalize as
# numbers (stripping trailing zeros, using scientific notation, etc). Force the
# Price cells to Text format before writing so the literal string is preserved,
# matching how the sheet already stores these columns as inline strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.431.87"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.727.47"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.96"
$ws.Range("E5").Value = "  -1.08%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4873"
$ws.Range("E7").Value = "  +1.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2617"
$ws.Range("E8").Value = "  -2.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06199"
$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.732.65"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07018"
$ws.Range("E11").Value = "  -2.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.45"
$ws.Range("E12").Value = "  -1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.556"
$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5982"
$ws.Range("E14").Value = "  -3.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.37"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.466.89"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007159"
$ws.Range("E19").Value = "  +2.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  -2.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.952.22"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.598"
$ws.Range("E23").Value = "  -3.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.180"
$ws.Range("E24").Value = "  -2.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.76"
$ws.Range("E25").Value = "  +1.55%  "

$ws.Range("E26").Value = "  -0.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.409"
$ws.Range("E27").Value = "  +0.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.93"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.722"
$ws.Range("E29").Value = "  -4.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.962"
$ws.Range("E30").Value = "  -0.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07985"
$ws.Range("E31").Value = "  -0.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.688"
$ws.Range("E32").Value = "  -0.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04523"
$ws.Range("E33").Value = "  -1.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9998"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.615"
$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  +0.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6237"
$ws.Range("E37").Value = "  -2.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9080"
$ws.Range("E38").Value = "  -1.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.983"
$ws.Range("E39").Value = "  -5.23%  "

$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01487"
$ws.Range("E42").Value = "  -1.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.29"
$ws.Range("E43").Value = "  -4.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.410"
$ws.Range("E44").Value = "  -3.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3867"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.687"
$ws.Range("E46").Value = "  -4.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1157"
$ws.Range("E47").Value = "  -2.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05362"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.26"
$ws.Range("E49").Value = "  -2.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.686"
$ws.Range("E50").Value = "  -2.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.256"
$ws.Range("E51").Value = "  -0.87%  "
